$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header cell F1 - same style as the other header cells (e.g. E1)
$ws.Range("E1").Copy()
$ws.Range("F1").PasteSpecial(-4122)
$ws.Range("F1").Value = "time_taken"

# time_taken values for rows 2..76 (data rows), written as plain text
$timeTaken = @(
    "2021-10-05 10:52:39.932192",
    "2021-10-05 10:52:39.932204",
    "2021-10-05 10:52:39.932207",
    "2021-10-05 10:52:39.932210",
    "2021-10-05 10:52:39.932212",
    "2021-10-05 10:52:39.932215",
    "2021-10-05 10:52:39.932218",
    "2021-10-05 10:52:39.932220",
    "2021-10-05 10:52:39.932223",
    "2021-10-05 10:52:39.932225",
    "2021-10-05 10:52:39.932228",
    "2021-10-05 10:52:39.932230",
    "2021-10-05 10:52:39.932233",
    "2021-10-05 10:52:39.932235",
    "2021-10-05 10:52:39.932238",
    "2021-10-05 10:52:39.932240",
    "2021-10-05 10:52:39.932243",
    "2021-10-05 10:52:39.932246",
    "2021-10-05 10:52:39.932248",
    "2021-10-05 10:52:39.932250",
    "2021-10-05 10:52:39.932253",
    "2021-10-05 10:52:39.932255",
    "2021-10-05 10:52:39.932258",
    "2021-10-05 10:52:39.932260",
    "2021-10-05 10:52:39.932263",
    "2021-10-05 10:52:39.932266",
    "2021-10-05 10:52:39.932268",
    "2021-10-05 10:52:39.932271",
    "2021-10-05 10:52:39.932273",
    "2021-10-05 10:52:39.932275",
    "2021-10-05 10:52:39.932278",
    "2021-10-05 10:52:39.932280",
    "2021-10-05 10:52:39.932283",
    "2021-10-05 10:52:39.932286",
    "2021-10-05 10:52:39.932288",
    "2021-10-05 10:52:39.932291",
    "2021-10-05 10:52:39.932293",
    "2021-10-05 10:52:39.932296",
    "2021-10-05 10:52:39.932298",
    "2021-10-05 10:52:39.932300",
    "2021-10-05 10:52:39.932303",
    "2021-10-05 10:52:39.932306",
    "2021-10-05 10:52:39.932308",
    "2021-10-05 10:52:39.932311",
    "2021-10-05 10:52:39.932313",
    "2021-10-05 10:52:39.932316",
    "2021-10-05 10:52:39.932318",
    "2021-10-05 10:52:39.932321",
    "2021-10-05 10:52:39.932323",
    "2021-10-05 10:52:39.932326",
    "2021-10-05 10:52:39.932328",
    "2021-10-05 10:52:39.932331",
    "2021-10-05 10:52:39.932334",
    "2021-10-05 10:52:39.932336",
    "2021-10-05 10:52:39.932339",
    "2021-10-05 10:52:39.932341",
    "2021-10-05 10:52:39.932344",
    "2021-10-05 10:52:39.932346",
    "2021-10-05 10:52:39.932349",
    "2021-10-05 10:52:39.932351",
    "2021-10-05 10:52:39.932353",
    "2021-10-05 10:52:39.932356",
    "2021-10-05 10:52:39.932359",
    "2021-10-05 10:52:39.932361",
    "2021-10-05 10:52:39.932365",
    "2021-10-05 10:52:39.932367",
    "2021-10-05 10:52:39.932370",
    "2021-10-05 10:52:39.932372",
    "2021-10-05 10:52:39.932375",
    "2021-10-05 10:52:39.932377",
    "2021-10-05 10:52:39.932380",
    "2021-10-05 10:52:39.932382",
    "2021-10-05 10:52:39.932385",
    "2021-10-05 10:52:39.932387",
    "2021-10-05 10:52:39.932390"
)

for ($i = 0; $i -lt $timeTaken.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 6).Value = $timeTaken[$i]
}
